$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column E (F2_Sales_Battle_2_Score) values for rows 2-13
$values = @(160, 150, 140, 130, 120, 110, 100, 90, 90, 90, 90, 90)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("E$row")
    $cell.Value = $values[$i]
}

# Match the style used by the neighboring D/F columns (red font) for E2:E13
$ws.Range("D2:D13").Copy() | Out-Null
$ws.Range("E2:E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the sheet view: remove frozen/scrolled topLeftCell and update the selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E2").Select() | Out-Null
